$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some cells in columns L/M are formatted as Text (numFmtId 49 = "@").
# Assigning a numeric .Value to such a cell stores it as a text string instead
# of a real number, so temporarily switch to General, set the number, then
# restore the original (Text) number format to keep the cell's style intact.
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Data corrections (raw input cells; dependent formulas recalc automatically) ---

# Row 274: new positive cases corrected down by 1
$ws.Range("C274").Value = 144

# Rows 282-292: intubation (E) / hors-SI (F) / SI (G) reclassification corrections
$ws.Range("E282").Value = 23
$ws.Range("F282").Value = 21
$ws.Range("G282").Value = 107

$ws.Range("E283").Value = 22
$ws.Range("F283").Value = 22
$ws.Range("G283").Value = 109

$ws.Range("E284").Value = 23
$ws.Range("F284").Value = 18
$ws.Range("G284").Value = 108

$ws.Range("E285").Value = 24
$ws.Range("F285").Value = 17
$ws.Range("G285").Value = 109

$ws.Range("E286").Value = 18
$ws.Range("F286").Value = 15
$ws.Range("G286").Value = 106

$ws.Range("E287").Value = 23
$ws.Range("F287").Value = 17
$ws.Range("G287").Value = 99

$ws.Range("E288").Value = 21
$ws.Range("F288").Value = 17
$ws.Range("G288").Value = 99

$ws.Range("E289").Value = 21
$ws.Range("F289").Value = 17
$ws.Range("G289").Value = 92

$ws.Range("E290").Value = 20
$ws.Range("G290").Value = 92

$ws.Range("E291").Value = 22
$ws.Range("G291").Value = 80

$ws.Range("E292").Value = 23
$ws.Range("G292").Value = 79

# Row 295: one additional in-hospital death recorded
Set-NumericValue $ws.Range("L295") 1

# Row 296 & 297: updated new-case counts
$ws.Range("C296").Value = 80
$ws.Range("C297").Value = 92

# Rows 298-300: fill in data that had previously been left blank
$ws.Range("C298").Value = 46
$ws.Range("E298").Value = 20
$ws.Range("F298").Value = 14
$ws.Range("G298").Value = 57
Set-NumericValue $ws.Range("L298") 0
Set-NumericValue $ws.Range("M298") 0

$ws.Range("C299").Value = 37
$ws.Range("E299").Value = 20
$ws.Range("F299").Value = 14
$ws.Range("G299").Value = 61
Set-NumericValue $ws.Range("L299") 0
Set-NumericValue $ws.Range("M299") 1

$ws.Range("C300").Value = 15
$ws.Range("E300").Value = 20
$ws.Range("F300").Value = 16
$ws.Range("G300").Value = 68
Set-NumericValue $ws.Range("L300") 0
Set-NumericValue $ws.Range("M300") 0

$excel.Calculate()

# --- View changes: scroll frozen pane back to top and select A2 ---
$ws.Range("A2").Select()
